# PlayerPerformance_3786.xlsx edit
# - Insert new "Player Info" sheet before "ODI Batting"
# - Rename MATCH_CARD_LINK -> MATCH_CODE in "ODI Batting" (col D) and "ODI Bowling" (col B),
#   replacing the howstat URL values with just the numeric match code
# - Drop the now-redundant empty INNING_NUMBER cells in "ODI Batting" (col B)
# - Append new "ODI Batting Extra" sheet at the end

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "ODI Batting" sheet (currently first sheet) - MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$lastRowBatting = $odiBatting.UsedRange.Rows.Count

$odiBatting.Range("D1").Value2 = "MATCH_CODE"

for ($r = 2; $r -le $lastRowBatting; $r++) {
    $linkCell = $odiBatting.Cells.Item($r, 4)
    $link = $linkCell.Value2
    if ($link -match "MatchCode=(\d+)") {
        $linkCell.NumberFormat = "@"
        $linkCell.Value2 = $matches[1]
    }

    $inningCell = $odiBatting.Cells.Item($r, 2)
    $inningVal = $inningCell.Value2
    if ($inningVal -eq "" -or $inningVal -eq $null) {
        $inningCell.ClearContents()
    }
}

# ---------------------------------------------------------------------------
# 2) "ODI Bowling" sheet - MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------------
$odiBowling = $wb.Worksheets.Item("ODI Bowling")
$lastRowBowling = $odiBowling.UsedRange.Rows.Count

$odiBowling.Range("B1").Value2 = "MATCH_CODE"

for ($r = 2; $r -le $lastRowBowling; $r++) {
    $linkCell = $odiBowling.Cells.Item($r, 2)
    $link = $linkCell.Value2
    if ($link -match "MatchCode=(\d+)") {
        $linkCell.NumberFormat = "@"
        $linkCell.Value2 = $matches[1]
    }
}

# ---------------------------------------------------------------------------
# 3) New "Player Info" sheet, inserted before "ODI Batting"
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($odiBatting)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value2 = "ID"
$playerInfo.Range("B1").Value2 = "NAME"
$playerInfo.Range("C1").Value2 = "BATTING_HAND"
$playerInfo.Range("D1").Value2 = "BOWL_STYLE"

$playerInfo.Range("A1:D1").Font.Bold = $true

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value2 = "3786"
$playerInfo.Range("B2").Value2 = "Umeshkumar Tilak Yadav"
$playerInfo.Range("C2").Value2 = "Right Handed"
$playerInfo.Range("D2").Value2 = "Right Arm Fast"

# ---------------------------------------------------------------------------
# 4) New "ODI Batting Extra" sheet, appended at the end
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra = $wb.Worksheets.Add($null, $lastSheet)
$battingExtra.Name = "ODI Batting Extra"

$battingExtra.Range("A1").Value2 = "MATCH_CODE"
$battingExtra.Range("B1").Value2 = "BATTING_POSITION"
$battingExtra.Range("C1").Value2 = "NUM_4"
$battingExtra.Range("D1").Value2 = "NUM_6"
$battingExtra.Range("E1").Value2 = "PERCENT_RUNS_OF_TOTAL"
$battingExtra.Range("F1").Value2 = "MAN_OF_MATCH"
$battingExtra.Range("A1:F1").Font.Bold = $true

$battingExtra.Columns.Item(1).NumberFormat = "@"
$battingExtra.Columns.Item(3).NumberFormat = "@"
$battingExtra.Columns.Item(4).NumberFormat = "@"
$battingExtra.Columns.Item(5).NumberFormat = "@"

$extraRows = @(
    @("3877", "",   "", "", "",      "NO"),
    @("3878", "",   "", "", "",      "NO"),
    @("3951", "",   "", "", "",      "NO"),
    @("3952", "10", "1", "0", "7.63%", "NO"),
    @("3953", "10", "", "", "",      "NO"),
    @("3954", "10", "", "", "",      "NO"),
    @("3955", "10", "", "", "",      "NO"),
    @("3974", "10", "", "", "",      "NO"),
    @("4034", "",   "", "", "",      "NO"),
    @("4038", "11", "", "", "",      "NO"),
    @("4051", "",   "", "", "",      "NO"),
    @("4052", "10", "0", "0", "",    "NO"),
    @("4053", "10", "", "", "",      "NO"),
    @("4056", "",   "", "", "",      "NO"),
    @("4057", "10", "", "", "",      "NO"),
    @("4074", "8",  "0", "0", "",    "NO"),
    @("4171", "",   "", "", "",      "NO"),
    @("4173", "10", "", "", "",      "NO"),
    @("4213", "",   "", "", "",      ""),
    @("4216", "",   "", "", "",      "")
)

$r = 2
foreach ($row in $extraRows) {
    $battingExtra.Cells.Item($r, 1).Value2 = $row[0]
    if ($row[1] -ne "") {
        $battingExtra.Cells.Item($r, 2).Value2 = [int]$row[1]
    }
    $battingExtra.Cells.Item($r, 3).Value2 = $row[2]
    $battingExtra.Cells.Item($r, 4).Value2 = $row[3]
    $battingExtra.Cells.Item($r, 5).Value2 = $row[4]
    if ($row[5] -ne "") {
        $battingExtra.Cells.Item($r, 6).Value2 = $row[5]
    }
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Leave selection on the first sheet, matching the original workbook state
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(1).Activate()
